# Update crypto price/volume table cells per scraper run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.608.82"
$ws.Range("E2").Value = "  +1.08%  "
$ws.Range("D3").Value = "2.989.83"
$ws.Range("E3").Value = "  +2.75%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'381.43"
$ws.Range("E5").Value = "  +4.59%  "
$ws.Range("D6").Value = "'106.27"
$ws.Range("E6").Value = "  +3.00%  "
$ws.Range("D7").Value = "'0.546"
$ws.Range("E7").Value = "  +1.05%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.600"
$ws.Range("E9").Value = "  +1.97%  "
$ws.Range("D10").Value = "'37.53"
$ws.Range("E10").Value = "  +1.65%  "
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("D12").Value = "'0.0846"
$ws.Range("E12").Value = "  +1.38%  "
$ws.Range("D13").Value = "'18.70"
$ws.Range("E13").Value = "  +1.65%  "
$ws.Range("D14").Value = "3.454.12"
$ws.Range("E14").Value = "  +2.42%  "
$ws.Range("D15").Value = "'7.53"
$ws.Range("E15").Value = "  +2.66%  "
$ws.Range("D16").Value = "2.982.04"
$ws.Range("E16").Value = "  +2.70%  "
$ws.Range("D17").Value = "'0.974"
$ws.Range("E17").Value = "  +2.52%  "
$ws.Range("D18").Value = "51.639.09"
$ws.Range("E18").Value = "  +1.20%  "
$ws.Range("D19").Value = "'3.41"
$ws.Range("E19").Value = "  +3.44%  "
$ws.Range("D20").Value = "'7.43"
$ws.Range("E20").Value = "  +2.65%  "
$ws.Range("D21").Value = "'13.06"
$ws.Range("E21").Value = "  +0.56%  "
$ws.Range("D22").Value = "0.0₃0964"
$ws.Range("E22").Value = "  +2.00%  "
$ws.Range("D23").Value = "'69.35"
$ws.Range("E23").Value = "  +1.85%  "
$ws.Range("D24").Value = "'263.85"
$ws.Range("E24").Value = "  +1.56%  "
$ws.Range("E25").Value = "  +4.41%  "
$ws.Range("E26").Value = "  -1.76%  "
$ws.Range("D27").Value = "'7.24"
$ws.Range("E27").Value = "  +18.62%  "
$ws.Range("D28").Value = "'7.50"
$ws.Range("E28").Value = "  +2.49%  "
$ws.Range("D30").Value = "'26.12"
$ws.Range("E30").Value = "  +0.97%  "
$ws.Range("E31").Value = "  +4.04%  "
$ws.Range("D32").Value = "'9.92"
$ws.Range("E32").Value = "  -0.18%  "
$ws.Range("D33").Value = "'35.09"
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("D34").Value = "'0.0463"
$ws.Range("E34").Value = "  +9.88%  "
$ws.Range("E35").Value = "  -2.05%  "
$ws.Range("D36").Value = "'51.38"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").Value = "'3.11"
$ws.Range("E38").Value = "  -0.30%  "
$ws.Range("D39").Value = "'17.51"
$ws.Range("E39").Value = "  +3.80%  "
$ws.Range("E40").Value = "  -5.84%  "
$ws.Range("D41").Value = "'1.86"
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("D42").Value = "'0.117"
$ws.Range("E42").Value = "  +2.85%  "
$ws.Range("D43").Value = "'123.68"
$ws.Range("E43").Value = "  +5.20%  "
$ws.Range("D44").Value = "'22.22"
$ws.Range("E44").Value = "  -0.88%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "'0.282"
$ws.Range("E45").Value = "  +20.30%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'2.08"
$ws.Range("E46").Value = "  -0.68%  "
$ws.Range("E47").Value = "  +5.05%  "
$ws.Range("D48").Value = "2.051.59"
$ws.Range("E48").Value = "  -0.48%  "
$ws.Range("D49").Value = "'3.26"
$ws.Range("E49").Value = "  +2.29%  "
$ws.Range("D50").Value = "'0.0355"
$ws.Range("E50").Value = "  +11.32%  "
$ws.Range("D51").Value = "'5.20"
$ws.Range("E51").Value = "  +3.62%  "
